$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row in the sheet
$lastRow = $ws.UsedRange.Rows.Count

# Swap the values of column C (codeforiati:group-code) and column D
# (codeforiati:group-name) for every row, including the header row.
for ($r = 1; $r -le $lastRow; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 3).Value2 = $dVal
    $ws.Cells.Item($r, 4).Value2 = $cVal
}
